$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the heading key in A1 with a custom label.
$ws.Range("A1").Value = "First grandchild"

# The "parent1.child3.*" heading/values that used to live in columns I:J
# now become the content of columns G:H (replacing the duplicated
# "parent1.child2.*" heading/values that used to be there).
$ws.Range("G1").Value = "parent1.child3.child31"
$ws.Range("H1").Value = "parent1.child3.child32"
$ws.Range("G2").ClearContents()
$ws.Range("H2").ClearContents()
$ws.Range("G3").Value = 30
$ws.Range("H3").Value = 31

# Columns I:L (the old "parent1.child3.*" / "parent2.child1.*" duplicate
# columns) are no longer used.
$ws.Range("I1:L3").ClearContents()
